# Update the "hvdc" worksheet: split the old "marginal_cost" column into a
# 3-term cost function (cost2, cost1, cost0), matching MATPOWER-style gencost
# fields. Existing marginal cost values move to the new "cost2" column, and
# two new columns "cost1" / "cost0" are added with value 0 for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hvdc")

# Rename header O1 from "marginal_cost" to "cost2"
$ws.Range("O1").Value = "cost2"

# Add the two new header cells
$ws.Range("P1").Value = "cost1"
$ws.Range("Q1").Value = "cost0"

# Populate the new columns for each data row (O column values are untouched)
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0

$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0

# Move the selection to match the saved view state
$ws.Range("Q2").Select()
